$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TextBox 6" holds the team-member list (right-aligned, red text) on slide 1.
$shp = $s.Shapes.Item("TextBox 6")
$tf = $shp.TextFrame.TextRange

# Current paragraphs (1-based character offsets into the shape's TextRange):
#   1: "Krishnaraj Palanychamy"   chars 1..22
#   2: "Krishnamurthy S"          chars 24..38
#   3: "Prabhakaran S"
#   4: "Pravin Kumar S"
#   5: "SR Abilash"
#   6: "Vishwanath Kannan"
#
# The edit swaps the names in paragraphs 1 and 2, and paragraph 2 ends up split
# into two runs ("Krishnaraj" / " Palanychamy") in the target deck.

# --- Paragraph 2: "Krishnamurthy S" -> "Krishnaraj Palanychamy" (as two runs) ---
# Replace the whole paragraph first (still a single run at this point).
$para2 = $tf.Characters(24, 15)
$para2.Text = "Krishnaraj Palanychamy"

# Re-assign just the "Krishnaraj" portion of that same range. Assigning text to a
# sub-range splits the run in two (new run for the assigned span, remainder keeps
# the original run's formatting) without altering any character formatting.
$run2a = $tf.Characters(24, 10)
$run2a.Text = "Krishnaraj"

# --- Paragraph 1: "Krishnaraj Palanychamy" -> "Krishnamurthy S" (single run) ---
$para1 = $tf.Characters(1, 22)
$para1.Text = "Krishnamurthy S"
